$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5607
$ws1.Range("F7").Value = 52
$ws1.Range("F10").Value = 5

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5607
$ws4.Range("F7").Value = 52
$ws4.Range("F11").Value = 5
